$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) cells are treated as text so values like "1.00" or "41.237.20"
# are preserved exactly (Excel would otherwise auto-convert numeric-looking strings to numbers).

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.237.20'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.09%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.172.64'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.25%  '

# Row 4
$ws.Range("E4").Value = '  -0.14%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.67'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.49%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.616'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.86%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '70.40'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -5.20%  '

# Row 8
$ws.Range("E8").Value = '  +0.02%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.580'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.70%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.39'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -8.44%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0930'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.69%  '

# Row 12
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.80'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -5.17%  '

# Row 13
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.101'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.82%  '

# Row 14
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.495.54'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.28%  '

# Row 15
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '13.94'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.75%  '

# Row 16
$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.807'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.91%  '

# Row 17
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.172.78'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.34%  '

# Row 18
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '40.998.11'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.21%  '

# Row 19
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0000102'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -8.12%  '

# Row 20
$ws.Range("B20").Value = 'Litecoin'
$ws.Range("C20").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.55'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.94%  '

# Row 21
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.95'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.71%  '

# Row 22
$ws.Range("B22").Value = 'InternetComputer(DFINITY)'
$ws.Range("C22").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.87'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.97%  '

# Row 23
$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '227.01'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.48%  '

# Row 24
$ws.Range("B24").Value = 'ImmutableX'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.96'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -7.19%  '

# Row 25
$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.28%  '

# Row 26
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.90'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.23%  '

# Row 27
$ws.Range("B27").Value = 'WEMIXToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.56'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.00%  '

# Row 28
$ws.Range("B28").Value = 'PancakeSwap'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.20'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.27%  '

# Row 29
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.19'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.63%  '

# Row 30
$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '167.89'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.53%  '

# Row 31
$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.01'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.35%  '

# Row 32
$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '30.54'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.04%  '

# Row 33
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0771'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.48%  '

# Row 34
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.18'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -9.26%  '

# Row 35
$ws.Range("B35").Value = 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.121'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.22%  '

# Row 36
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.103'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -10.05%  '

# Row 37
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.12'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.31%  '

# Row 38
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0286'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.92%  '

# Row 39
$ws.Range("B39").Value = 'Celestia'
$ws.Range("C39").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '12.03'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.04%  '

# Row 40
$ws.Range("B40").Value = 'LidoDAOToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.09'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.24%  '

# Row 41
$ws.Range("B41").Value = 'THORChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.45'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.62%  '

# Row 42
$ws.Range("B42").Value = 'MultiversX'
$ws.Range("C42").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '60.22'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -10.00%  '

# Row 43
$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.191'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.90%  '

# Row 44
$ws.Range("B44").Value = 'Cronos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0978'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.37%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.31'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.68%  '

# Row 46
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '98.21'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.22%  '

# Row 47
$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.09'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.86%  '

# Row 48
$ws.Range("B48").Value = 'TrustWalletToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.14'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.41%  '

# Row 49
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.22'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -9.35%  '

# Row 50
$ws.Range("B50").Value = 'HuobiToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.63'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.90%  '

# Row 51
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.374.22'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.27%  '
